# Regenerate instance to have positive average demands during the last periods.
$wb = $excel.ActiveWorkbook

# --- Productdata sheet: Gamma (column G) values for first four products ---
$wsProd = $wb.Worksheets.Item("Productdata")
$wsProd.Range("G2").Value = 49
$wsProd.Range("G3").Value = 21
$wsProd.Range("G4").Value = 35
$wsProd.Range("G5").Value = 70

# Guard against the engine's round-trip quirk where an empty shared-string
# cell (t="s" with no <v>, i.e. column H's blank placeholder) picks up a
# stray value on save. Re-asserting blanks keeps those cells truly empty.
for ($r = 2; $r -le 11; $r++) {
    $wsProd.Range("H$r").Value = ""
}

# --- ForecastedAverageDemand sheet: fill in demand for the last three buckets ---
$wsAvg = $wb.Worksheets.Item("ForecastedAverageDemand")
foreach ($r in 9..11) {
    $wsAvg.Range("B$r").Value = 70
    $wsAvg.Range("C$r").Value = 30
    $wsAvg.Range("D$r").Value = 50
    $wsAvg.Range("E$r").Value = 100
}

# --- ForcastedStandardDeviation sheet: matching standard deviations ---
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStd.Range("B9").Value = 7.166424999999998
$wsStd.Range("C9").Value = 3.071324999999999
$wsStd.Range("D9").Value = 5.118874999999999
$wsStd.Range("E9").Value = 10.23775

$wsStd.Range("B10").Value = 8.1997825
$wsStd.Range("C10").Value = 3.5141925
$wsStd.Range("D10").Value = 5.856987499999999
$wsStd.Range("E10").Value = 11.713975

$wsStd.Range("B11").Value = 9.129804249999998
$wsStd.Range("C11").Value = 3.912773249999999
$wsStd.Range("D11").Value = 6.521288749999998
$wsStd.Range("E11").Value = 13.0425775
